# Refresh Leve profit-calculation inputs/outputs (currentAveragePrice*,
# LevePrice*, LeveProfit*) on each job sheet with latest market data,
# per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 21
$ws.Cells.Item(21, 8).Value = 0  # H21
$ws.Cells.Item(21, 9).Value = 0  # I21
$ws.Cells.Item(21, 11).Value = 0  # K21
$ws.Cells.Item(21, 13).ClearContents()  # M21

# row 23
$ws.Cells.Item(23, 8).Value = 0  # H23
$ws.Cells.Item(23, 9).Value = 0  # I23
$ws.Cells.Item(23, 11).Value = 0  # K23
$ws.Cells.Item(23, 13).ClearContents()  # M23

# row 28
$ws.Cells.Item(28, 8).Value = 1444.421  # H28
$ws.Cells.Item(28, 9).Value = 555.5294  # I28
$ws.Cells.Item(28, 10).Value = 9000  # J28
$ws.Cells.Item(28, 11).Value = 555.5294  # K28
$ws.Cells.Item(28, 12).Value = 9000  # L28
$ws.Cells.Item(28, 13).Value = -70.52940000000001  # M28
$ws.Cells.Item(28, 14).Value = -9970  # N28

# row 33
$ws.Cells.Item(33, 8).Value = 157.125  # H33
$ws.Cells.Item(33, 9).Value = 157.52174  # I33
$ws.Cells.Item(33, 10).Value = 148  # J33
$ws.Cells.Item(33, 11).Value = 157.52174  # K33
$ws.Cells.Item(33, 12).Value = 148  # L33
$ws.Cells.Item(33, 13).Value = 71.47826000000001  # M33
$ws.Cells.Item(33, 14).Value = -606  # N33

# row 62
$ws.Cells.Item(62, 8).Value = 3000.1  # H62
$ws.Cells.Item(62, 9).Value = 2570.4285  # I62
$ws.Cells.Item(62, 10).Value = 4002.6667  # J62
$ws.Cells.Item(62, 11).Value = 2570.4285  # K62
$ws.Cells.Item(62, 12).Value = 4002.6667  # L62
$ws.Cells.Item(62, 13).Value = -1946.4285  # M62
$ws.Cells.Item(62, 14).Value = -5250.6667  # N62

# row 65
$ws.Cells.Item(65, 8).Value = 3000.1  # H65
$ws.Cells.Item(65, 9).Value = 2570.4285  # I65
$ws.Cells.Item(65, 10).Value = 4002.6667  # J65
$ws.Cells.Item(65, 11).Value = 12852.1425  # K65
$ws.Cells.Item(65, 12).Value = 20013.3335  # L65
$ws.Cells.Item(65, 13).Value = -9732.1425  # M65
$ws.Cells.Item(65, 14).Value = -26253.3335  # N65

# row 137
$ws.Cells.Item(137, 8).Value = 1356248.2  # H137
$ws.Cells.Item(137, 9).Value = 2174620.8  # I137
$ws.Cells.Item(137, 11).Value = 6523862.399999999  # K137
$ws.Cells.Item(137, 13).Value = -6521312.399999999  # M137

# row 141
$ws.Cells.Item(141, 8).Value = 3840.5881  # H141
$ws.Cells.Item(141, 9).Value = 2032.6666  # I141
$ws.Cells.Item(141, 11).Value = 6097.9998  # K141
$ws.Cells.Item(141, 13).Value = -917.9997999999996  # M141

$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Cells.Item(45, 8).Value = 17544.371  # H45
$ws.Cells.Item(45, 10).Value = 6009  # J45
$ws.Cells.Item(45, 12).Value = 6009  # L45
$ws.Cells.Item(45, 14).Value = -6763  # N45

# row 61
$ws.Cells.Item(61, 8).Value = 2247.532  # H61
$ws.Cells.Item(61, 9).Value = 1741.2903  # I61
$ws.Cells.Item(61, 10).Value = 3228.375  # J61
$ws.Cells.Item(61, 11).Value = 1741.2903  # K61
$ws.Cells.Item(61, 12).Value = 3228.375  # L61
$ws.Cells.Item(61, 13).Value = -1529.2903  # M61
$ws.Cells.Item(61, 14).Value = -3652.375  # N61

# row 88
$ws.Cells.Item(88, 8).Value = 3172.2856  # H88
$ws.Cells.Item(88, 9).Value = 1592.1666  # I88
$ws.Cells.Item(88, 11).Value = 1592.1666  # K88
$ws.Cells.Item(88, 13).Value = -1186.1666  # M88

# row 91
$ws.Cells.Item(91, 8).Value = 3172.2856  # H91
$ws.Cells.Item(91, 9).Value = 1592.1666  # I91
$ws.Cells.Item(91, 11).Value = 1592.1666  # K91
$ws.Cells.Item(91, 13).Value = -188.1666  # M91

# row 132
$ws.Cells.Item(132, 8).Value = 1784.762  # H132
$ws.Cells.Item(132, 9).Value = 952.21875  # I132
$ws.Cells.Item(132, 11).Value = 2856.65625  # K132
$ws.Cells.Item(132, 13).Value = -326.65625  # M132

# row 136
$ws.Cells.Item(136, 8).Value = 2247.532  # H136
$ws.Cells.Item(136, 9).Value = 1741.2903  # I136
$ws.Cells.Item(136, 10).Value = 3228.375  # J136
$ws.Cells.Item(136, 11).Value = 5223.8709  # K136
$ws.Cells.Item(136, 12).Value = 9685.125  # L136
$ws.Cells.Item(136, 13).Value = -2673.8709  # M136
$ws.Cells.Item(136, 14).Value = -14785.125  # N136

$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Cells.Item(86, 8).Value = 2285.762  # H86
$ws.Cells.Item(86, 9).Value = 2142  # I86
$ws.Cells.Item(86, 11).Value = 2142  # K86
$ws.Cells.Item(86, 13).Value = -1019  # M86

# row 89
$ws.Cells.Item(89, 8).Value = 2285.762  # H89
$ws.Cells.Item(89, 9).Value = 2142  # I89
$ws.Cells.Item(89, 11).Value = 10710  # K89
$ws.Cells.Item(89, 13).Value = -5094  # M89

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Cells.Item(31, 8).Value = 3331.3333  # H31
$ws.Cells.Item(31, 9).Value = 1903.2273  # I31
$ws.Cells.Item(31, 11).Value = 1903.2273  # K31
$ws.Cells.Item(31, 13).Value = -1608.2273  # M31

# row 34
$ws.Cells.Item(34, 8).Value = 3331.3333  # H34
$ws.Cells.Item(34, 9).Value = 1903.2273  # I34
$ws.Cells.Item(34, 11).Value = 1903.2273  # K34
$ws.Cells.Item(34, 13).Value = -1701.2273  # M34

# row 58
$ws.Cells.Item(58, 8).Value = 2851.074  # H58
$ws.Cells.Item(58, 9).Value = 2200.3076  # I58
$ws.Cells.Item(58, 10).Value = 3455.3572  # J58
$ws.Cells.Item(58, 11).Value = 2200.3076  # K58
$ws.Cells.Item(58, 12).Value = 3455.3572  # L58
$ws.Cells.Item(58, 13).Value = -1997.3076  # M58
$ws.Cells.Item(58, 14).Value = -3861.3572  # N58

# row 103
$ws.Cells.Item(103, 8).Value = 12456.286  # H103
$ws.Cells.Item(103, 9).Value = 12456.286  # I103
$ws.Cells.Item(103, 11).Value = 12456.286  # K103
$ws.Cells.Item(103, 13).Value = -11284.286  # M103

# row 132
$ws.Cells.Item(132, 8).Value = 4156.846  # H132
$ws.Cells.Item(132, 9).Value = 4425.6665  # I132
$ws.Cells.Item(132, 11).Value = 13276.9995  # K132
$ws.Cells.Item(132, 13).Value = -10746.9995  # M132

# row 136
$ws.Cells.Item(136, 8).Value = 2851.074  # H136
$ws.Cells.Item(136, 9).Value = 2200.3076  # I136
$ws.Cells.Item(136, 10).Value = 3455.3572  # J136
$ws.Cells.Item(136, 11).Value = 6600.9228  # K136
$ws.Cells.Item(136, 12).Value = 10366.0716  # L136
$ws.Cells.Item(136, 13).Value = -4050.9228  # M136
$ws.Cells.Item(136, 14).Value = -15466.0716  # N136

$ws = $wb.Worksheets.Item("CUL")
# row 39
$ws.Cells.Item(39, 8).Value = 6553.2144  # H39
$ws.Cells.Item(39, 10).Value = 6553.2144  # J39
$ws.Cells.Item(39, 12).Value = 19659.6432  # L39
$ws.Cells.Item(39, 14).Value = -20247.6432  # N39

# row 116
$ws.Cells.Item(116, 8).Value = 114543.375  # H116
$ws.Cells.Item(116, 9).Value = 180390.6  # I116
$ws.Cells.Item(116, 11).Value = 541171.8  # K116
$ws.Cells.Item(116, 13).Value = -537729.8  # M116

# row 117
$ws.Cells.Item(117, 8).Value = 4333  # H117
$ws.Cells.Item(117, 10).Value = 4999.5  # J117
$ws.Cells.Item(117, 12).Value = 14998.5  # L117
$ws.Cells.Item(117, 14).Value = -21882.5  # N117

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Cells.Item(80, 8).Value = 40001740  # H80
$ws.Cells.Item(80, 9).Value = 125002070  # I80
$ws.Cells.Item(80, 10).Value = 1585.2354  # J80
$ws.Cells.Item(80, 11).Value = 125002070  # K80
$ws.Cells.Item(80, 12).Value = 1585.2354  # L80
$ws.Cells.Item(80, 13).Value = -125001072  # M80
$ws.Cells.Item(80, 14).Value = -3581.2354  # N80

# row 83
$ws.Cells.Item(83, 8).Value = 40001740  # H83
$ws.Cells.Item(83, 9).Value = 125002070  # I83
$ws.Cells.Item(83, 10).Value = 1585.2354  # J83
$ws.Cells.Item(83, 11).Value = 625010350  # K83
$ws.Cells.Item(83, 12).Value = 7926.177  # L83
$ws.Cells.Item(83, 13).Value = -625005358  # M83
$ws.Cells.Item(83, 14).Value = -17910.177  # N83

# row 113
$ws.Cells.Item(113, 8).Value = 5581.1665  # H113
$ws.Cells.Item(113, 9).Value = 4847.8  # I113
$ws.Cells.Item(113, 11).Value = 4847.8  # K113
$ws.Cells.Item(113, 13).Value = -2677.8  # M113

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Cells.Item(22, 8).Value = 750.375  # H22
$ws.Cells.Item(22, 9).Value = 700.6  # I22
$ws.Cells.Item(22, 11).Value = 700.6  # K22
$ws.Cells.Item(22, 13).Value = -405.6  # M22

# row 27
$ws.Cells.Item(27, 8).Value = 750.375  # H27
$ws.Cells.Item(27, 9).Value = 700.6  # I27
$ws.Cells.Item(27, 11).Value = 700.6  # K27
$ws.Cells.Item(27, 13).Value = -593.6  # M27

# row 40
$ws.Cells.Item(40, 8).Value = 6670.5625  # H40
$ws.Cells.Item(40, 10).Value = 6949.7334  # J40
$ws.Cells.Item(40, 12).Value = 6949.7334  # L40
$ws.Cells.Item(40, 14).Value = -7221.7334  # N40

# row 56
$ws.Cells.Item(56, 8).Value = 35136.25  # H56
$ws.Cells.Item(56, 9).Value = 25272.5  # I56
$ws.Cells.Item(56, 11).Value = 25272.5  # K56
$ws.Cells.Item(56, 13).Value = -24581.5  # M56

# row 68
$ws.Cells.Item(68, 8).Value = 3956.8333  # H68
$ws.Cells.Item(68, 10).Value = 3999  # J68
$ws.Cells.Item(68, 12).Value = 3999  # L68
$ws.Cells.Item(68, 14).Value = -5497  # N68

# row 71
$ws.Cells.Item(71, 8).Value = 3956.8333  # H71
$ws.Cells.Item(71, 10).Value = 3999  # J71
$ws.Cells.Item(71, 12).Value = 19995  # L71
$ws.Cells.Item(71, 14).Value = -27483  # N71

# row 82
$ws.Cells.Item(82, 8).Value = 729.2727  # H82
$ws.Cells.Item(82, 9).Value = 947.5  # I82
$ws.Cells.Item(82, 10).Value = 467.4  # J82
$ws.Cells.Item(82, 11).Value = 947.5  # K82
$ws.Cells.Item(82, 12).Value = 467.4  # L82
$ws.Cells.Item(82, 13).Value = -586.5  # M82
$ws.Cells.Item(82, 14).Value = -1189.4  # N82

# row 85
$ws.Cells.Item(85, 8).Value = 729.2727  # H85
$ws.Cells.Item(85, 9).Value = 947.5  # I85
$ws.Cells.Item(85, 10).Value = 467.4  # J85
$ws.Cells.Item(85, 11).Value = 947.5  # K85
$ws.Cells.Item(85, 12).Value = 467.4  # L85
$ws.Cells.Item(85, 13).Value = 300.5  # M85
$ws.Cells.Item(85, 14).Value = -2963.4  # N85

# row 94
$ws.Cells.Item(94, 8).Value = 0  # H94
$ws.Cells.Item(94, 10).Value = 0  # J94
$ws.Cells.Item(94, 12).Value = 0  # L94
$ws.Cells.Item(94, 14).ClearContents()  # N94

# row 122
$ws.Cells.Item(122, 8).Value = 2930.85  # H122
$ws.Cells.Item(122, 9).Value = 2844.2727  # I122
$ws.Cells.Item(122, 10).Value = 3036.6667  # J122
$ws.Cells.Item(122, 11).Value = 8532.8181  # K122
$ws.Cells.Item(122, 12).Value = 9110.000100000001  # L122
$ws.Cells.Item(122, 13).Value = -6082.8181  # M122
$ws.Cells.Item(122, 14).Value = -14010.0001  # N122

# row 132
$ws.Cells.Item(132, 8).Value = 6128.0386  # H132
$ws.Cells.Item(132, 9).Value = 3507.2727  # I132
$ws.Cells.Item(132, 10).Value = 8049.933  # J132
$ws.Cells.Item(132, 11).Value = 10521.8181  # K132
$ws.Cells.Item(132, 12).Value = 24149.799  # L132
$ws.Cells.Item(132, 13).Value = -7991.8181  # M132
$ws.Cells.Item(132, 14).Value = -29209.799  # N132

# row 136
$ws.Cells.Item(136, 8).Value = 5728.4165  # H136
$ws.Cells.Item(136, 9).Value = 5592.75  # I136
$ws.Cells.Item(136, 11).Value = 16778.25  # K136
$ws.Cells.Item(136, 13).Value = -14228.25  # M136

$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Cells.Item(81, 8).Value = 4271.875  # H81
$ws.Cells.Item(81, 9).Value = 2664.2856  # I81
$ws.Cells.Item(81, 10).Value = 5522.222  # J81
$ws.Cells.Item(81, 11).Value = 5328.5712  # K81
$ws.Cells.Item(81, 12).Value = 11044.444  # L81
$ws.Cells.Item(81, 13).Value = -4267.5712  # M81
$ws.Cells.Item(81, 14).Value = -13166.444  # N81

# row 84
$ws.Cells.Item(84, 8).Value = 4271.875  # H84
$ws.Cells.Item(84, 9).Value = 2664.2856  # I84
$ws.Cells.Item(84, 10).Value = 5522.222  # J84
$ws.Cells.Item(84, 11).Value = 26642.856  # K84
$ws.Cells.Item(84, 12).Value = 55222.22  # L84
$ws.Cells.Item(84, 13).Value = -21338.856  # M84
$ws.Cells.Item(84, 14).Value = -65830.22  # N84

# row 136
$ws.Cells.Item(136, 8).Value = 66671024  # H136
$ws.Cells.Item(136, 9).Value = 76924160  # I136
$ws.Cells.Item(136, 11).Value = 230772480  # K136
$ws.Cells.Item(136, 13).Value = -230769930  # M136
